# Applies numeric updates to H:N columns across several sheets
# as part of the scheduled profit-data refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 2846.8
$ws.Range("I141").Value = 1801.0834
$ws.Range("K141").Value = 5403.2502
$ws.Range("M141").Value = -223.2502000000004

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 999.7273
$ws.Range("I2").Value = 999.7273
$ws.Range("K2").Value = 999.7273
$ws.Range("M2").Value = -886.7273

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 240.35715
$ws.Range("I5").Value = 212
$ws.Range("J5").Value = 344.33334
$ws.Range("K5").Value = 212
$ws.Range("L5").Value = 344.33334
$ws.Range("M5").Value = -100
$ws.Range("N5").Value = -568.33334

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H31").Value = 31591.666
$ws.Range("I31").Value = 31591.666
$ws.Range("K31").Value = 31591.666
$ws.Range("M31").Value = -31297.666

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 999.7273
$ws.Range("I116").Value = 999.7273
$ws.Range("K116").Value = 999.7273
$ws.Range("M116").Value = 1294.2727

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 999.7273
$ws.Range("I3").Value = 999.7273
$ws.Range("K3").Value = 999.7273
$ws.Range("M3").Value = -885.7273

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 240.35715
$ws.Range("I4").Value = 212
$ws.Range("J4").Value = 344.33334
$ws.Range("K4").Value = 212
$ws.Range("L4").Value = 344.33334
$ws.Range("M4").Value = -97
$ws.Range("N4").Value = -574.33334

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3436.2856
$ws.Range("I20").Value = 3814.1667
$ws.Range("J20").Value = 1169
$ws.Range("K20").Value = 3814.1667
$ws.Range("L20").Value = 1169
$ws.Range("M20").Value = -3567.1667
$ws.Range("N20").Value = -1663

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H26").Value = 43867.5
$ws.Range("I26").Value = 43867.5
$ws.Range("K26").Value = 43867.5
$ws.Range("M26").Value = -43575.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H27").Value = 0
$ws.Range("J27").Value = 0
$ws.Range("L27").Value = 0
$ws.Range("N27").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H36").Value = 991.53845
$ws.Range("I36").Value = 815
$ws.Range("K36").Value = 815
$ws.Range("M36").Value = -281

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H96").Value = 19713.375
$ws.Range("I96").Value = 19672.428
$ws.Range("K96").Value = 19672.428
$ws.Range("M96").Value = -16926.428

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H97").Value = 1616.3334
$ws.Range("I97").Value = 1616.3334
$ws.Range("K97").Value = 1616.3334
$ws.Range("M97").Value = -625.3334

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 947.375
$ws.Range("I99").Value = 947.375
$ws.Range("K99").Value = 947.375
$ws.Range("M99").Value = 550.625

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H101").Value = 0
$ws.Range("J101").Value = 0
$ws.Range("L101").Value = 0
$ws.Range("N101").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 5415
$ws.Range("I134").Value = 1246.9231
$ws.Range("K134").Value = 3740.7693
$ws.Range("M134").Value = -1205.7693

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 1634.3334
$ws.Range("I62").Value = 905
$ws.Range("J62").Value = 1999
$ws.Range("K62").Value = 905
$ws.Range("L62").Value = 1999
$ws.Range("M62").Value = -281
$ws.Range("N62").Value = -3247

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H65").Value = 1634.3334
$ws.Range("I65").Value = 905
$ws.Range("J65").Value = 1999
$ws.Range("K65").Value = 4525
$ws.Range("L65").Value = 9995
$ws.Range("M65").Value = -1405
$ws.Range("N65").Value = -16235

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H68").Value = 98765
$ws.Range("J68").Value = 98765
$ws.Range("L68").Value = 98765
$ws.Range("N68").Value = -100263

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H71").Value = 98765
$ws.Range("J71").Value = 98765
$ws.Range("L71").Value = 296295
$ws.Range("N71").Value = -303783

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H128").Value = 428567
$ws.Range("I128").Value = 428567
$ws.Range("K128").Value = 1285701
$ws.Range("M128").Value = -1280721

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 4613.7144
$ws.Range("I80").Value = 4539.2
$ws.Range("K80").Value = 4539.2
$ws.Range("M80").Value = -3541.2

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 4613.7144
$ws.Range("I83").Value = 4539.2
$ws.Range("K83").Value = 22696
$ws.Range("M83").Value = -17704

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 1374.3334
$ws.Range("I122").Value = 670.6
$ws.Range("K122").Value = 2011.8
$ws.Range("M122").Value = 438.1999999999998

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H20").Value = 503500
$ws.Range("J20").Value = 503500
$ws.Range("L20").Value = 503500
$ws.Range("N20").Value = -503952

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1539.4375
$ws.Range("I22").Value = 1233.1428
$ws.Range("J22").Value = 1777.6666
$ws.Range("K22").Value = 1233.1428
$ws.Range("L22").Value = 1777.6666
$ws.Range("M22").Value = -938.1428000000001
$ws.Range("N22").Value = -2367.6666

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 1539.4375
$ws.Range("I27").Value = 1233.1428
$ws.Range("J27").Value = 1777.6666
$ws.Range("K27").Value = 1233.1428
$ws.Range("L27").Value = 1777.6666
$ws.Range("M27").Value = -1126.1428
$ws.Range("N27").Value = -1991.6666

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 3153.4614
$ws.Range("I82").Value = 1070.7142
$ws.Range("K82").Value = 1070.7142
$ws.Range("M82").Value = -709.7141999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H85").Value = 3153.4614
$ws.Range("I85").Value = 1070.7142
$ws.Range("K85").Value = 1070.7142
$ws.Range("M85").Value = 177.2858000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H104").Value = 18839.334
$ws.Range("J104").Value = 18839.334
$ws.Range("L104").Value = 18839.334
$ws.Range("N104").Value = -25827.334

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 1200
$ws.Range("I122").Value = 1100
$ws.Range("J122").Value = 1300
$ws.Range("K122").Value = 3300
$ws.Range("L122").Value = 3900
$ws.Range("M122").Value = -850
$ws.Range("N122").Value = -8800

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 0
$ws.Range("I132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("M132").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H30").Value = 15000
$ws.Range("J30").Value = 15000
$ws.Range("L30").Value = 15000
$ws.Range("N30").Value = -15214

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 8739
$ws.Range("I62").Value = 3217.3333
$ws.Range("K62").Value = 3217.3333
$ws.Range("M62").Value = -2593.3333

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H65").Value = 8739
$ws.Range("I65").Value = 3217.3333
$ws.Range("K65").Value = 16086.6665
$ws.Range("M65").Value = -12966.6665

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H104").Value = 20130.334
$ws.Range("J104").Value = 20130.334
$ws.Range("L104").Value = 20130.334
$ws.Range("N104").Value = -27118.334

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2982.0715
$ws.Range("I122").Value = 1745.9
$ws.Range("J122").Value = 6072.5
$ws.Range("K122").Value = 5237.700000000001
$ws.Range("L122").Value = 18217.5
$ws.Range("M122").Value = -2787.700000000001
$ws.Range("N122").Value = -23117.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 915.3333
$ws.Range("I132").Value = 915.3333
$ws.Range("K132").Value = 2745.9999
$ws.Range("M132").Value = -215.9998999999998

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()
